$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country names (shared string reorder) in column A ---
$ws.Range("A71").Value = "Azerbaiyan"
$ws.Range("A72").Value = "Eslovenia"
$ws.Range("A73").Value = "Banglades"
$ws.Range("A100").Value = "Guinea"
$ws.Range("A101").Value = "Malta"
$ws.Range("A102").Value = "Bolivia"
$ws.Range("A103").Value = "Jordania"
$ws.Range("A104").Value = "Taiwan"
$ws.Range("A105").Value = "Reunion"
$ws.Range("A106").Value = "Nigeria"
$ws.Range("A107").Value = "San Marino"
$ws.Range("A149").Value = "Cabo Verde"
$ws.Range("A150").Value = "Polinesia Francesa"
$ws.Range("A151").Value = "Uganda"
$ws.Range("A152").Value = "Islas Caimanes"
$ws.Range("A153").Value = "San Martin (Parte Holandesa)"
$ws.Range("A154").Value = "Bahamas"
$ws.Range("A156").Value = "Zambia"

# --- Update footer text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 15:52"

# --- Update numeric data cells ---
$ws.Range("B4").Value = 614643
$ws.Range("C4").Value = 757
$ws.Range("E4").Value = 549711
$ws.Range("G4").Value = 65
$ws.Range("H4").Value = 26112
$ws.Range("B8").Value = 132321
$ws.Range("C8").Value = 111
$ws.Range("E8").Value = 56219
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = 3502
$ws.Range("B9").Value = 98476
$ws.Range("C9").Value = 4603
$ws.Range("E9").Value = 85264
$ws.Range("G9").Value = 761
$ws.Range("H9").Value = 12868
$ws.Range("B32").Value = 6740
$ws.Range("C32").Value = 117
$ws.Range("E32").Value = 6563
$ws.Range("G32").Value = 6
$ws.Range("H32").Value = 145
$ws.Range("F43").Value = 128
$ws.Range("E51").Value = 2865
$ws.Range("G51").Value = 8
$ws.Range("H51").Value = 72
$ws.Range("D54").Value = 596
$ws.Range("E54").Value = 1738
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 109
$ws.Range("B62").Value = 1727
$ws.Range("C62").Value = 7
$ws.Range("D62").Value = 1077
$ws.Range("E62").Value = 642
$ws.Range("B71").Value = 1253
$ws.Range("C71").Value = 56
$ws.Range("D71").Value = 404
$ws.Range("E71").Value = 836
$ws.Range("F71").Value = 24
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 13
$ws.Range("B72").Value = 1248
$ws.Range("C72").Value = 28
$ws.Range("D72").Value = 165
$ws.Range("E72").Value = 1022
$ws.Range("F72").Value = 34
$ws.Range("G72").Value = 5
$ws.Range("H72").Value = 61
$ws.Range("B73").Value = 1231
$ws.Range("C73").Value = 219
$ws.Range("D73").Value = 49
$ws.Range("E73").Value = 1132
$ws.Range("F73").Value = 1
$ws.Range("G73").Value = 4
$ws.Range("H73").Value = 50
$ws.Range("B100").Value = 404
$ws.Range("C100").Value = 41
$ws.Range("D100").Value = 31
$ws.Range("E100").Value = 372
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 1
$ws.Range("B101").Value = 399
$ws.Range("C101").Value = 6
$ws.Range("D101").Value = 44
$ws.Range("E101").Value = 352
$ws.Range("F101").Value = 4
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 3
$ws.Range("C102").Value = 43
$ws.Range("D102").Value = 7
$ws.Range("E102").Value = 362
$ws.Range("F102").Value = 3
$ws.Range("H102").Value = 28
$ws.Range("B103").Value = 397
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 235
$ws.Range("E103").Value = 155
$ws.Range("F103").Value = 5
$ws.Range("H103").Value = 7
$ws.Range("B104").Value = 395
$ws.Range("C104").Value = 2
$ws.Range("D104").Value = 137
$ws.Range("E104").Value = 252
$ws.Range("F104").Value = 0
$ws.Range("H104").Value = 6
$ws.Range("B105").Value = 391
$ws.Range("D105").Value = 40
$ws.Range("E105").Value = 351
$ws.Range("F105").Value = 3
$ws.Range("H105").Value = 0
$ws.Range("B106").Value = 373
$ws.Range("D106").Value = 99
$ws.Range("E106").Value = 263
$ws.Range("F106").Value = 2
$ws.Range("H106").Value = 11
$ws.Range("B107").Value = 372
$ws.Range("D107").Value = 53
$ws.Range("E107").Value = 283
$ws.Range("F107").Value = 15
$ws.Range("H107").Value = 36
$ws.Range("B116").Value = 254
$ws.Range("C116").Value = 13
$ws.Range("D116").Value = 21
$ws.Range("E116").Value = 212
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 21
$ws.Range("B133").Value = 110
$ws.Range("C133").Value = 2
$ws.Range("D133").Value = 29
$ws.Range("E133").Value = 81
$ws.Range("B149").Value = 56
$ws.Range("C149").Value = 45
$ws.Range("D149").Value = 1
$ws.Range("E149").Value = 54
$ws.Range("F149").Value = 0
$ws.Range("H149").Value = 1
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 55
$ws.Range("F150").Value = 1
$ws.Range("B151").Value = 55
$ws.Range("D151").Value = 12
$ws.Range("E151").Value = 43
$ws.Range("F151").Value = 0
$ws.Range("H151").Value = 0
$ws.Range("B152").Value = 54
$ws.Range("D152").Value = 6
$ws.Range("E152").Value = 47
$ws.Range("F152").Value = 3
$ws.Range("H152").Value = 1
$ws.Range("B153").Value = 52
$ws.Range("D153").Value = 5
$ws.Range("E153").Value = 38
$ws.Range("F153").Value = 2
$ws.Range("H153").Value = 9
$ws.Range("B154").Value = 49
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 6
$ws.Range("E154").Value = 35
$ws.Range("H154").Value = 8
$ws.Range("B155").Value = 48
$ws.Range("C155").Value = 1
$ws.Range("E155").Value = 34
$ws.Range("B156").Value = 48
$ws.Range("C156").Value = 3
$ws.Range("D156").Value = 30
$ws.Range("E156").Value = 16
$ws.Range("F156").Value = 1
$ws.Range("H156").Value = 2
$ws.Range("D157").Value = 15
$ws.Range("E157").Value = 30
